# Applies the "Updated symbol list" commit: refreshed price/volume figures
# and a 3-way reorder of the BKEXToken/CEJI/KickToken rows (41-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    # Force the cell to Text format so numeric-looking strings (including
    # ones with significant trailing zeros) are preserved verbatim, exactly
    # like the source inlineStr cells.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# --- Column D (Price) updates -------------------------------------------
Set-TextValue "D2"  "243.25"
Set-TextValue "D4"  "5.287"
Set-TextValue "D6"  "6.476"
Set-TextValue "D7"  "3.332"
Set-TextValue "D8"  "0.8086"
Set-TextValue "D9"  "0.8800"
Set-TextValue "D10" "0.1382"
Set-TextValue "D11" "0.07281"
Set-TextValue "D12" "0.03088"
Set-TextValue "D13" "0.03057"
Set-TextValue "D14" "0.09326"
Set-TextValue "D15" "3.862"
Set-TextValue "D16" "0.001539"
Set-TextValue "D19" "0.005897"
Set-TextValue "D22" "0.00008700"
Set-TextValue "D23" "3.578"
Set-TextValue "D25" "0.3183"
Set-TextValue "D40" "0.03771"
Set-TextValue "D44" "0.007103"
Set-TextValue "D45" "0.00005480"

# --- E18: drop the stray "Worstin24h" suffix -----------------------------
$ws.Range("E18").Value = "17OneONE"

# --- Rows 41-43: BKEXToken / CEJI / KickToken reshuffle -------------------
# New order: row41 = KickToken, row42 = BKEXToken, row43 = CEJI
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006410"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1053"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002560"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
